# Fix the typo in the header of column D ("departemen" -> "departement")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "departement"

# Update the "licence" column (E) values for every student row from 5 to 3
$ws.Range("E2:E57").Value = 3

# Leave the selection on the cell that was last edited, as Excel would
$ws.Range("D1").Select()
